# coldbrook-feeding.xlsx: add a "comment row" (row 4) under the header row,
# replacing the old cell-comment (NoteBook) annotations with plain,
# always-visible guidance cells. Matches commit "add comment row to all
# templates, fix movement unit tests".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the legacy cell comments that used to carry this guidance
#    (A3:J3). Deleting every comment drops comments1.xml, the VML
#    drawing relationship and the <legacyDrawing> element on save.
# ------------------------------------------------------------------
$commentCells = @("A3", "B3", "C3", "D3", "E3", "F3", "G3", "H3", "I3", "J3")
foreach ($addr in $commentCells) {
    $cmt = $ws.Range($addr).Comment
    if ($cmt -ne $null) {
        $cmt.Delete()
    }
}

# ------------------------------------------------------------------
# 2) Write the new guidance row (row 4) as real cell values.
# ------------------------------------------------------------------
$ws.Range("A4").Value = "Name of the Tank"
$ws.Range("B4").Value = "Must match feeding method in database. Eg Automatic"
$ws.Range("C4").Value = "Optional. Frequency of feeding. Eg. Daily, 1/minute."
$ws.Range("D4").Value = "Optional"
$ws.Range("E4").Value = "Feed type must match feed code in database"
$ws.Range("F4").Value = "Size of feed, numeric. Eg. 0.5, 1, 3"
$ws.Range("G4").Value = "Feed type must match feed code in database"
$ws.Range("H4").Value = "Size of feed, numeric. Eg. 0.5, 1, 3"
$ws.Range("I4").Value = "Feed type must match feed code in database"
$ws.Range("J4").Value = "Size of feed, numeric. Eg. 0.5, 1, 3"
$ws.Range("K4").Value = "Use same format for columns as feed types 1-3"

# ------------------------------------------------------------------
# 3) Format the new row: small grey Tahoma font, top-aligned wrapped
#    text, and a thin box border around each cell (A4 skips its own
#    left edge since it sits at the sheet boundary).
# ------------------------------------------------------------------
$commentRow = $ws.Range("A4:K4")
$commentRow.Font.Name = "Tahoma"
$commentRow.Font.Size = 10
$commentRow.Font.Color = 8421504   # RGB(128,128,128) = FF808080
$commentRow.VerticalAlignment = -4160   # xlTop
$commentRow.WrapText = $true

$ws.Range("B4:K4").Borders.LineStyle = 1
$ws.Range("B4:K4").Borders.Weight = 2

$a4 = $ws.Range("A4")
$a4.Borders.LineStyle = 1
$a4.Borders.Weight = 2
$a4.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> xlLineStyleNone

$ws.Rows.Item(4).RowHeight = 51

# ------------------------------------------------------------------
# 4) Widen the new comments column (K) to fit the extra text.
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 31.43

Write-Output "done"
